$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.970149755477905
$ws.Range("B1").Value = 2.629374027252197
$ws.Range("C1").Value = 2.277615308761597
$ws.Range("D1").Value = 2.404961347579956
$ws.Range("E1").Value = 3.087450504302979
